# Adds a bold space + the YouTube URL run right after the
# "URL to Your Coding Assignment Video:" run, before the existing
# line-break run, reproducing the target diff's two new <w:r> runs.

$d = $word.ActiveDocument

$marker = "URL to Your Coding Assignment Video:"

# Locate the paragraph that contains the marker text.
$searchRng = $d.Content
$found = $searchRng.Find.Execute($marker, $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '$marker' in the document."
}

# Re-resolve the *document-level* paragraph (not the Find hit's own,
# possibly-clipped sub-range) so its .Range covers the whole paragraph,
# including the trailing <w:br/> run and paragraph mark.
$paraIndex = $searchRng.Paragraphs.Item(1).Index
$para = $d.Paragraphs.Item($paraIndex)
$paraRng = $para.Range

# Pull the paragraph's real OOXML (with its original attributes intact)
# out of the WordOpenXML package payload.
$openXml = $paraRng.WordOpenXML
if (-not ($openXml -match '(?s)<w:body>(.*?)</w:body>')) {
    throw "Could not locate <w:body> in WordOpenXML."
}
$body = $matches[1]
if (-not ($body -match '(?s)^(<w:p\b.*?</w:p>)')) {
    throw "Could not locate the target <w:p> in WordOpenXML."
}
$paraXml = $matches[1]

# The two new runs to splice in right after the marker's run closes.
$insertion = '<w:r><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr>' + `
             '<w:t xml:space="preserve"> </w:t></w:r>' + `
             '<w:r><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr>' + `
             '<w:t>https://youtu.be/DUK_hz6I1eI</w:t></w:r>'

$afterMarkerRun = "$marker</w:t></w:r>"
$insertAt = $paraXml.IndexOf($afterMarkerRun)
if ($insertAt -lt 0) {
    throw "Could not find the marker run boundary in the paragraph XML."
}
$insertAt = $insertAt + $afterMarkerRun.Length

$newParaXml = $paraXml.Substring(0, $insertAt) + $insertion + `
              $paraXml.Substring($insertAt)

# InsertXML replaces the contents of the exact range it is called on,
# so calling it on the whole-paragraph range swaps the paragraph's runs
# for our modified copy (original run + 2 new runs + original <w:br/>).
$paraRng.InsertXML($newParaXml)
